$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume updates from the "Updated symbol list" data refresh.
# Cells are stored as text in the workbook, so force Text format before
# assigning numeric-looking / percent-looking strings, then restore the
# original (default) style so no stray formatting is introduced.
$cellUpdates = @{
    'D2' = '288.96'
    'E2' = '-4.70%'
    'D3' = '30.22'
    'E3' = '-7.45%'
    'D4' = '4.913'
    'E4' = '-2.92%'
    'D5' = '0.07231'
    'E5' = '-6.38%'
    'D6' = '1.795'
    'E6' = '-14.68%'
    'D7' = '7.590'
    'E7' = '-3.95%'
    'D8' = '3.718'
    'E8' = '-1.98%'
    'D9' = '0.9018'
    'E9' = '-2.36%'
    'D10' = '0.1666'
    'E10' = '-5.74%'
    'D11' = '0.07853'
    'E11' = '-1.63%'
    'D12' = '0.07966'
    'E12' = '-7.10%'
    'D13' = '0.03044'
    'E13' = '-0.98%'
    'D14' = '0.1003'
    'E14' = '0.57%'
    'D15' = '0.001493'
    'E15' = '-1.78%'
    'D16' = '0.005775'
    'E16' = '2.24%'
    'D17' = '3.472'
    'E17' = '-0.10%'
    'D18' = '2.076'
    'E18' = '-3.73%'
    'D19' = '0.3318'
    'E19' = '-0.60%'
    'E20' = '-1.78%'
    'D21' = '3.984'
    'E21' = '-9.32%'
    'D22' = '0.2173'
    'E22' = '10.11%'
    'D23' = '0.04502'
    'E23' = '-0.87%'
    'D24' = '0.001218'
    'E24' = '-0.92%'
    'D25' = '0.004414'
    'E25' = '6.33%'
    'D26' = '0.0001300'
    'E26' = '4.13%'
    'D27' = '0.0003395'
    'E27' = '-95.47%'
    'D39' = '0.01565'
    'E39' = '-9.66%'
    'D40' = '0.04328'
    'E40' = '-8.14%'
    'D41' = '0.007323'
    'E41' = '-1.90%'
    'D42' = '0.01006'
    'D43' = '0.1307'
    'E43' = '-3.96%'
    'D44' = '0.002006'
    'E44' = '-13.83%'
    'D45' = '0.009045'
    'E45' = '-14.64%'
    'D46' = '0.00005882'
    'E46' = '-4.35%'
    'D47' = '0.00000000750'
    'E47' = '0.11%'
    'E48' = '102.85%'
    'D49' = '0.002899'
    'E49' = '-3.26%'
    'D50' = '0.00002100'
    'E50' = '0.11%'
    'D51' = '0.0002000'
    'E51' = '0.11%'
}

foreach ($addr in $cellUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$addr]
    $cell.Style = "Normal"
}
